$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for the new "Save" column (H), matching the existing header
# formatting (bold, bordered, centered) used by the other header cells.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Fill H2:H57 with 1 when the "sum" column (G) is greater than 8, else 0.
for ($r = 2; $r -le 57; $r++) {
    $g = $ws.Cells.Item($r, 7).Value2
    if ($g -gt 8) {
        $ws.Cells.Item($r, 8).Value = 1
    } else {
        $ws.Cells.Item($r, 8).Value = 0
    }
}
